$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast-error values (naive component forecaster bug fix)
$ws.Range("J25").Value = -0.4802901972718701
$ws.Range("K25").Value = 0.7862757877985587
$ws.Range("I26").Value = -0.360290197271857
$ws.Range("J26").Value = 0.9062757877985718
$ws.Range("H27").Value = -0.3802901972718758
$ws.Range("I27").Value = 0.8662757877985712
$ws.Range("G28").Value = -0.360290197271857
$ws.Range("H28").Value = 0.8062757877985718
$ws.Range("F29").Value = -0.3202901972718593
$ws.Range("G29").Value = 0.9462757877985695
$ws.Range("H29").Value = 0.240008545046976
$ws.Range("I29").Value = -1.369651846369791
$ws.Range("J29").Value = 2.65850917732945
$ws.Range("K29").Value = -0.5180856546984671
$ws.Range("E30").Value = -0.2602901972718712
$ws.Range("F30").Value = 1.006275787798558
$ws.Range("G30").Value = 0.3000085450469641
$ws.Range("H30").Value = -1.309651846369803
$ws.Range("I30").Value = 2.698509177329442
$ws.Range("J30").Value = -0.478085654698475
$ws.Range("D31").Value = -0.2602901972718712
$ws.Range("E31").Value = 0.9662757877985655
$ws.Range("F31").Value = 0.2900085450469732
$ws.Range("G31").Value = -1.319651846369794
$ws.Range("H31").Value = 2.708509177329447
$ws.Range("I31").Value = -0.46808565469847
$ws.Range("C32").Value = 0.6397098027281429
$ws.Range("D32").Value = 1.306275787798572
$ws.Range("E32").Value = 0.1000085450469783
$ws.Range("F32").Value = -1.409651846369789
$ws.Range("G32").Value = 2.618509177329452
$ws.Range("H32").Value = -0.5580856546984648
$ws.Range("B33").Value = -0.4102901982718521
$ws.Range("C33").Value = 1.106275787798566
$ws.Range("D33").Value = 0.200008545046984
$ws.Range("E33").Value = -1.409651846369783
$ws.Range("F33").Value = 2.618509177329458
$ws.Range("G33").Value = -0.5580856546984592
$ws.Range("H33").Value = 0.3315486822015572
$ws.Range("I33").Value = -1.944079640195553
$ws.Range("J33").Value = 0.40000558765135
$ws.Range("K33").Value = -0.2233623182505511
$ws.Range("B34").Value = 1.206275787798589
$ws.Range("C34").Value = 0.200008545046984
$ws.Range("D34").Value = -1.409651846369783
$ws.Range("E34").Value = 2.418509177329455
$ws.Range("F34").Value = -0.758085654698462
$ws.Range("G34").Value = 0.1315486822015543
$ws.Range("H34").Value = -2.144079640195556
$ws.Range("I34").Value = 0.2000055876513471
$ws.Range("J34").Value = -0.5233623182505625
$ws.Range("B35").Value = 0.5000085450469953
$ws.Range("C35").Value = -0.7096518463697949
$ws.Range("D35").Value = 2.518509177329435
$ws.Range("E35").Value = -0.758085654698462
$ws.Range("F35").Value = 0.03154868220154583
$ws.Range("G35").Value = -2.244079640195565
$ws.Range("H35").Value = 0.1000055876513386
$ws.Range("I35").Value = -0.5233623182505625
$ws.Range("B36").Value = -1.009651846369792
$ws.Range("C36").Value = 2.718509177329438
$ws.Range("D36").Value = -1.358085654698456
$ws.Range("E36").Value = -0.6684513177984428
$ws.Range("F36").Value = -2.144079640195556
$ws.Range("G36").Value = 0.5000055876513301
$ws.Range("H36").Value = -0.123362318250571
$ws.Range("B37").Value = 1.418509177329455
$ws.Range("C37").Value = -1.658085654698468
$ws.Range("D37").Value = -0.06845131779844849
$ws.Range("E37").Value = -2.244079640195565
$ws.Range("F37").Value = 0.1000055876513386
$ws.Range("G37").Value = -0.5233623182505625
$ws.Range("H37").Value = 4.987355777631223
$ws.Range("I37").Value = 0.02253521981180029
$ws.Range("J37").Value = -0.1231996868497021
$ws.Range("K37").Value = -0.2473129037311424
$ws.Range("B38").Value = -0.4580856546984791
$ws.Range("C38").Value = 0.2315486822015344
$ws.Range("D38").Value = -2.944079640195553
$ws.Range("E38").Value = -0.59999441234865
$ws.Range("F38").Value = -0.5233623182505625
$ws.Range("G38").Value = 5.037355777631234
$ws.Range("H38").Value = 0.1225352198118088
$ws.Range("I38").Value = -0.02319968684969362
$ws.Range("J38").Value = -0.1473129037311623
$ws.Range("B39").Value = -0.06845131779842006
$ws.Range("C39").Value = -2.544079640195562
$ws.Range("D39").Value = -0.3999944123486614
$ws.Range("E39").Value = -0.3233623182505739
$ws.Range("F39").Value = 5.087355777631231
$ws.Range("G39").Value = 0.1225352198118088
$ws.Range("H39").Value = -0.02319968684969362
$ws.Range("I39").Value = -0.1473129037311339
$ws.Range("B40").Value = -1.844079640195545
$ws.Range("C40").Value = 0.1000055876513386
$ws.Range("D40").Value = -0.3233623182505739
$ws.Range("E40").Value = 5.087355777631231
$ws.Range("F40").Value = 0.2725352198118145
$ws.Range("G40").Value = 0.07680031315028649
$ws.Range("H40").Value = -0.04731290373115382
$ws.Range("B41").Value = 0.7000055876513613
$ws.Range("C41").Value = -0.8233623182505312
$ws.Range("D41").Value = 5.387355777631214
$ws.Range("E41").Value = 0.4225352198117918
$ws.Range("F41").Value = 0.07680031315028649
$ws.Range("G41").Value = -0.04731290373115382
$ws.Range("H41").Value = 0.7268782796002142
$ws.Range("I41").Value = -0.7346752961032477
$ws.Range("J41").Value = -0.1343669380733701
$ws.Range("K41").Value = -0.1160091367467971
$ws.Range("B42").Value = -1.623362318250543
$ws.Range("C42").Value = 5.087355777631231
$ws.Range("D42").Value = 0.7225352198118031
$ws.Range("E42").Value = 0.1768003131503093
$ws.Range("F42").Value = -0.2473129037311424
$ws.Range("G42").Value = 0.7268782796002284
$ws.Range("H42").Value = -0.5346752961032732
$ws.Range("I42").Value = -0.3343669380733587
$ws.Range("J42").Value = -0.3160091367468141
$ws.Range("B43").Value = 4.387355777631228
$ws.Range("C43").Value = 0.3225352198118117
$ws.Range("D43").Value = 0.1768003131503093
$ws.Range("E43").Value = -0.2473129037311424
$ws.Range("F43").Value = 0.7268782796002284
$ws.Range("G43").Value = -0.7346752961032761
$ws.Range("H43").Value = -0.1343669380733701
$ws.Range("I43").Value = -0.1160091367467971
$ws.Range("B44").Value = 2.722535219811803
$ws.Range("C44").Value = -0.2231996868496964
$ws.Range("D44").Value = 0.2526870962688576
$ws.Range("E44").Value = 0.8268782796002512
$ws.Range("F44").Value = -0.8346752961032562
$ws.Range("G44").Value = -0.2343669380733502
$ws.Range("H44").Value = -0.2160091367467772
$ws.Range("B45").Value = -0.2231996868496964
$ws.Range("C45").Value = -0.04731290373115382
$ws.Range("D45").Value = 0.9268782796002313
$ws.Range("E45").Value = -0.8346752961032562
$ws.Range("F45").Value = -0.2343669380733502
$ws.Range("G45").Value = -0.2160091367467772
$ws.Range("H45").Value = 0.6423454266220716
$ws.Range("I45").Value = -0.5919400257838703
$ws.Range("B46").Value = -0.1473129037311339
$ws.Range("C46").Value = 0.6268782796002483
$ws.Range("D46").Value = -0.8346752961032562
$ws.Range("E46").Value = -0.1343669380733701
$ws.Range("F46").Value = -0.1160091367467971
$ws.Range("G46").Value = 0.7423454266220944
$ws.Range("H46").Value = -0.3919400257838674
$ws.Range("B47").Value = 0.3268782796002512
$ws.Range("C47").Value = -1.134675296103239
$ws.Range("D47").Value = -0.3343669380733587
$ws.Range("E47").Value = -0.2160091367467772
$ws.Range("F47").Value = 0.7423454266220944
$ws.Range("G47").Value = -0.3919400257838674
$ws.Range("B48").Value = -0.8346752961032562
$ws.Range("C48").Value = -0.2343669380733502
$ws.Range("D48").Value = -0.1160091367467971
$ws.Range("E48").Value = 0.7023454266220881
$ws.Range("F48").Value = -0.4319400257838737
$ws.Range("B49").Value = -0.3343669380733587
$ws.Range("C49").Value = -0.3160091367467857
$ws.Range("D49").Value = 0.5423454266220915
$ws.Range("E49").Value = -0.4919400257838902
$ws.Range("B50").Value = -0.3160091367467857
$ws.Range("C50").Value = 0.442345426622083
$ws.Range("D50").Value = -0.5919400257838703
$ws.Range("B51").Value = 0.3423454266220887
$ws.Range("C51").Value = -0.7919400257838731
$ws.Range("B52").Value = -0.2919400257838873

# Clear trailing cells that are no longer part of the diagonal band
$ws.Range("J45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("H47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("B53").ClearContents()
